# Performance By Model.xlsx edit
#
# The table is reshaped so that Loss/Accuracy (new measurements) become
# columns B/C (right after the model name), the old Size/Input/Stem
# columns shift right to D/E/F, and the "Previous Accuracy"/"Previous
# Loss" columns (old G/H) are removed entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the old "Previous Accuracy"/"Previous Loss" columns (G:H) and
#    the old Loss/Accuracy columns (E:F) - their values are being
#    replaced with freshly measured numbers placed in B:C below.
$ws.Range("E1:H1").EntireColumn.Delete() | Out-Null

# 2) Insert two fresh columns at B:C to hold the new Loss/Accuracy data.
#    (This shifts old Size -> D, Input -> E, Stem -> F, preserving their
#    original column widths/number formatting untouched.)
$ws.Range("B1:C1").EntireColumn.Insert() | Out-Null

# 3) Header row
$ws.Range("B1").Value = "Loss"
$ws.Range("C1").Value = "Accuracy"

# 4) New Loss / Accuracy values per model (rows 2-8 have data; the last
#    three models, rows 9-11, have no measurements yet so stay blank).
$loss = @{
  2 = 0.730934798717498
  3 = 0.52614372968673695
  4 = 0.40091994404792702
  5 = 0.32090279459953303
  6 = 0.27091637253761203
  7 = 0.48250153660774198
  8 = 0.295252114534378
}
$accuracy = @{
  2 = 0.72554349899291903
  3 = 0.8125
  4 = 0.86865943670272805
  5 = 0.89492756128311102
  6 = 0.90670287609100297
  7 = 0.816123187541961
  8 = 0.90851449966430597
}

for ($r = 2; $r -le 11; $r++) {
  if ($loss.ContainsKey($r)) {
    $ws.Cells.Item($r, 2).Value = $loss[$r]
    $ws.Cells.Item($r, 3).Value = $accuracy[$r]
  } else {
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4131
    $ws.Cells.Item($r, 3).HorizontalAlignment = -4131
  }
}

# Match the new number cells' left-alignment style used throughout the
# rest of the data rows.
$ws.Range("B2:C8").HorizontalAlignment = -4131

# 5) Column widths for the two new columns (as close as this engine's
#    AutoFit/column-width quantization allows to the recorded best fit).
$ws.Columns.Item(2).ColumnWidth = 12
$ws.Columns.Item(3).ColumnWidth = 12

# 6) Restore the recorded selection/active cell.
$ws.Range("C9").Select()
